$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 565
$wsExhibit.Range("F4").Value = 1233
$wsExhibit.Range("F5").Value = 1061
$wsExhibit.Range("F6").Value = 13993
$wsExhibit.Range("F7").Value = 15254
$wsExhibit.Range("F9").Value = 32
$wsExhibit.Range("F23").Value = 5915
$wsExhibit.Range("F24").Value = 955
$wsExhibit.Range("F25").Value = 1075
$wsExhibit.Range("F26").Value = 5497
$wsExhibit.Range("F28").Value = 134
$wsExhibit.Range("F29").Value = 80
$wsExhibit.Range("F30").Value = 414

# Sheet "全部类型" (sheet4) updates to column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 565
$wsAll.Range("F5").Value = 1233
$wsAll.Range("F6").Value = 1061
$wsAll.Range("F7").Value = 13993
$wsAll.Range("F8").Value = 15254
$wsAll.Range("F10").Value = 32
$wsAll.Range("F25").Value = 5915
$wsAll.Range("F26").Value = 955
$wsAll.Range("F27").Value = 1075
$wsAll.Range("F28").Value = 5497
$wsAll.Range("F30").Value = 134
$wsAll.Range("F31").Value = 80
$wsAll.Range("F32").Value = 414

$wb.Save()
